# Fill in the half-year report card comments and grades for
# Math (מתמטיקה), English (אנגלית) and History of Israel (תולדות ישראל).
# Each subject block in the document looks like:
#   <subject heading paragraph>
#   <comment paragraph>   (empty run inside a vMerge "restart" cell)
#   ...
#   <"ciyon:" (grade label) paragraph>
#   <grade paragraph>     (empty run, to receive the numeric grade)
#
# We locate every block by matching the subject heading text, then fill
# the paragraph right after it (the comment box) and the paragraph right
# after the following "ciyon:" paragraph (the grade box).

$d = $word.ActiveDocument
$nl = [char]10

$updates = @{
    "מתמטיקה"      = @{ Comment = ("במחצית למדנו משוואות ב2 נעלמים, פיתחנו כמה שיטות לבעיה זו,בנוסף התעסקנו בבעיות תנועה וזמן ולמדנו איך להתמודד מול זאת" + $nl + "היה לנו הספקים מעולים!" + $nl + "אילה את ילדה מקסימה, שיהיה לך הרבה הצלחה בהמשך! "); Grade = "98" }
    "אנגלית"       = @{ Comment = ("במחצית זאת התמקדנו על הבנה חזקה של הטקסטים ולמדנו את השיטות להבנת הנקרא, חזרנו על שאלות חוזרות ופיתחנו שיטות קלות לפיתרתן." + $nl + "אילה את ילדה נהדרת, הרבה הצלחה!"); Grade = "100" }
    "תולדות ישראל" = @{ Comment = ("במחצית זאת למדנו על גדולי ישראל בכל מיני יבשות, על המצב של היהודים בתקופות שלטון שונות," + $nl + "אילה הרבה הצלחה!"); Grade = "97" }
}

$total = $d.Paragraphs.Count

foreach ($subject in $updates.Keys) {
    $info = $updates[$subject]
    $subjectIdx = -1

    for ($i = 1; $i -le $total; $i++) {
        $ptext = $d.Paragraphs.Item($i).Range.Text
        $trimmed = $ptext.TrimEnd([char]13, [char]7, [char]32)
        if ($trimmed -eq $subject) {
            $subjectIdx = $i
            break
        }
    }

    if ($subjectIdx -eq -1) {
        Write-Host "WARNING: subject not found:" $subject
        continue
    }

    # Comment paragraph immediately follows the subject heading paragraph.
    $commentPara = $d.Paragraphs.Item($subjectIdx + 1)
    $commentPara.Range.Text = $info.Comment

    # Find the "ciyon:" (grade label) paragraph that follows this subject,
    # then fill the paragraph right after it with the grade.
    $ciyonIdx = -1
    for ($i = $subjectIdx + 1; $i -le $total; $i++) {
        $ptext = $d.Paragraphs.Item($i).Range.Text
        if ($ptext -match "ציון:") {
            $ciyonIdx = $i
            break
        }
    }

    if ($ciyonIdx -eq -1) {
        Write-Host "WARNING: ciyon label not found for subject:" $subject
        continue
    }

    $gradePara = $d.Paragraphs.Item($ciyonIdx + 1)
    $gradePara.Range.Text = $info.Grade
}

Write-Host "Done."
